$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D, shifting D:J -> E:K
$ws.Range("D1:D6").EntireColumn.Insert()

# Populate the new column D (PR - polynomial regression)
$ws.Range("D1").Value = "PR"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

$ws.Range("D2").Value = "GridSearchCV"
$ws.Range("D3").Value = 12
$ws.Range("D4").Value = "38m, 29s"
$ws.Range("D5").Value = "23m, 55s"
$ws.Range("D6").Value = "15s"

# Update the remaining cells with final tuned/computed results
$ws.Range("B4").Value = "39s"
$ws.Range("E4").Value = "41s"
$ws.Range("F4").Value = "14m, 10s"
$ws.Range("G4").Value = "19m, 46s"
$ws.Range("H4").Value = "11m, 51s"
$ws.Range("I4").Value = "1m, 39s"
$ws.Range("J4").Value = "35m, 37s"
$ws.Range("K4").Value = "1h, 37m"

$ws.Range("E5").Value = "2m, 13s"
$ws.Range("F5").Value = "54m, 22s"
$ws.Range("G5").Value = "1h, 14m"
$ws.Range("H5").Value = "44m, 3s"
$ws.Range("I5").Value = "5m, 47s"
$ws.Range("J5").Value = "1h, 42m"
$ws.Range("K5").Value = "6h, 26m"

$ws.Range("B6").Value = "2m, 12s"
$ws.Range("E6").Value = "2s"
$ws.Range("F6").Value = "11s"
$ws.Range("G6").Value = "6s"
$ws.Range("H6").Value = "1m, 38s"
$ws.Range("I6").Value = "2s"
$ws.Range("J6").Value = "36m, 22s"
$ws.Range("K6").Value = "16s"
